# Table of new Price (D) / Volume(1h) (E) text values for the refreshed
# coin rows (GitHub Actions symbol-list update, 2023-02-06 15:23 UTC).
# Keys are row numbers on Sheet1; D/E hold the literal text to write.
$updates = @{
    2 = @{ D="327.31"; E="-0.89%" }
    3 = @{ D="43.85"; E="0.65%" }
    4 = @{ D="5.544"; E="-0.78%" }
    5 = @{ D="0.08029"; E="-2.17%" }
    6 = @{ D="1.900"; E="0.19%" }
    7 = @{ D="4.268"; E="-2.33%" }
    8 = @{ D="0.9464"; E="0.33%" }
    9 = @{ D="2.540"; E="-9.83%" }
    10 = @{ D="0.1180"; E="-0.98%" }
    11 = @{ D="0.1844"; E="-3.83%" }
    12 = @{ D="0.09631"; E="-2.62%" }
    13 = @{ D="0.04400"; E="1.35%" }
    14 = @{ D="0.1067"; E="-0.11%" }
    15 = @{ D="0.001286"; E="0.58%" }
    16 = @{ D="0.005987"; E="0.37%" }
    17 = @{ D="3.405"; E="-3.57%" }
    18 = @{ D="0.3447"; E="-2.54%" }
    19 = @{ D="10.58"; E="21.06%" }
    20 = @{ D="0.1379"; E="0.72%" }
    21 = @{ D="0.2507"; E="0.47%" }
    22 = @{ D="0.04181"; E="-4.67%" }
    23 = @{ D="0.001247"; E="0.47%" }
    24 = @{ D="0.004285"; E="-1.29%" }
    25 = @{ D="0.0001262"; E="2.17%" }
    26 = @{ D="0.0003994"; E="-0.24%" }
    38 = @{ D="0.02645"; E="-5.24%" }
    39 = @{ D="0.05500"; E="-3.96%" }
    40 = @{ D="0.007588"; E="-4.50%" }
    41 = @{ E="-1.98%" }
    42 = @{ D="0.007984"; E="-18.10%" }
    43 = @{ D="0.002003"; E="-4.81%" }
    44 = @{ D="0.008801"; E="-12.50%" }
    45 = @{ D="0.00006895"; E="-9.21%" }
    46 = @{ D="0.00000000751"; E="-0.30%" }
    47 = @{ D="0.003522"; E="2.13%" }
    48 = @{ D="0.002273"; E="-0.24%" }
    49 = @{ D="0.00002102"; E="-0.30%" }
    50 = @{ D="0.0002002"; E="-0.30%" }
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row].Keys) {
        $cell = $ws.Range("$col$row")
        # Leading apostrophe forces text storage (prevents Excel from
        # reinterpreting '327.31' as a number or '-0.89%' as a percentage,
        # which would drop significant trailing zeros like "1.900").
        $cell.Value = "'$($updates[$row][$col])"
        # Reset the style Excel auto-assigns for quote-prefixed text so the
        # cell keeps its original (default/no explicit style) formatting.
        $cell.Style = "Normal"
    }
}
